$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Kebabas13!"
$ws.Range("C3").Value = "Pieštukas17!"
$ws.Range("C4").Value = "Lėktuvas23!"
$ws.Range("C5").Value = "Zebras33!"
$ws.Range("C6").Value = "ABC!!23??a"
$ws.Range("C7").Value = "ABC??12aa"
$ws.Range("C8").Value = "ErelisLėktuve!i2"
$ws.Range("C9").Value = "Sparnuotas?O2?"
$ws.Range("C10").Value = "125521ABCa??"
$ws.Range("C11").Value = "ABCCABa1??"

$ws.Range("C7").Select()
